$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.233.84"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.860.22"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +1.32%  "
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("D6").Value = "'311.90"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "'0.4797"
$ws.Range("E7").Value = "  +2.33%  "
$ws.Range("D8").Value = "'0.3724"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").Value = "'0.9365"
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("D11").Value = "'20.27"
$ws.Range("E11").Value = "  +3.61%  "
$ws.Range("D12").Value = "'0.07874"
$ws.Range("E12").Value = "  +2.46%  "
$ws.Range("D13").Value = "1.851.17"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'5.422"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "'6.532"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("D16").Value = "'90.26"
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").Value = "'0.000008743"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "27.254.27"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "'14.70"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "'5.102"
$ws.Range("E22").Value = "  +1.75%  "
$ws.Range("D23").Value = "'10.66"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'1.948"
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("D25").Value = "'153.83"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'18.49"
$ws.Range("D27").Value = "'1.996"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'115.61"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'4.970"
$ws.Range("E29").Value = "  +1.92%  "
$ws.Range("D30").Value = "'0.08887"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "'3.346"
$ws.Range("E31").Value = "  +4.16%  "
$ws.Range("D32").Value = "'1.184"
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D33").Value = "'4.586"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").Value = "'0.7419"
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "'2.687"
$ws.Range("E35").Value = "  -3.30%  "
$ws.Range("D36").Value = "'1.124"
$ws.Range("E36").Value = "  +3.79%  "
$ws.Range("D37").Value = "'0.02026"
$ws.Range("E37").Value = "  +4.47%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05258"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5330"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'7.113"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1526"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'8.337"
$ws.Range("E42").Value = "  +2.40%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'10.60"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.4786"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.020"
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'102.61"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.635"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'66.33"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06075"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'0.8981"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'36.74"
$ws.Range("E51").Value = "  +1.28%  "
